$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.966.64"
$ws.Range("E2").Value = "  -4.74%  "

$ws.Range("D3").Value = "2.494.72"
$ws.Range("E3").Value = "  -3.25%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.46"
$ws.Range("E5").Value = "  -2.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.72"
$ws.Range("E6").Value = "  -6.98%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.34%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.573"
$ws.Range("E8").Value = "  -3.20%  "

$ws.Range("D9").Value = "2.528.68"
$ws.Range("E9").Value = "  -2.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0996"
$ws.Range("E10").Value = "  -4.18%  "

$ws.Range("E11").Value = "  -2.72%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.44"
$ws.Range("E12").Value = "  -0.97%  "

$ws.Range("E13").Value = "  -4.05%  "

$ws.Range("D14").Value = "2.937.72"
$ws.Range("E14").Value = "  -3.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.74"
$ws.Range("E15").Value = "  -6.60%  "

$ws.Range("D16").Value = "58.867.75"
$ws.Range("E16").Value = "  -4.80%  "

$ws.Range("E17").Value = "  -3.56%  "

$ws.Range("D18").Value = "2.510.44"
$ws.Range("E18").Value = "  -2.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.26"
$ws.Range("E19").Value = "  -2.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.27"
$ws.Range("E20").Value = "  -5.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.88"
$ws.Range("E21").Value = "  -4.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.996"
$ws.Range("E22").Value = "  -0.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.74"
$ws.Range("E23").Value = "  -4.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.78"
$ws.Range("E24").Value = "  -2.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.438"
$ws.Range("E25").Value = "  -10.69%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.162"
$ws.Range("E26").Value = "  -3.31%  "

$ws.Range("D27").Value = "2.610.59"
$ws.Range("E27").Value = "  -3.17%  "

$ws.Range("E28").Value = "  -0.65%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.74"
$ws.Range("E29").Value = "  -4.65%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.76"
$ws.Range("E30").Value = "  -6.12%  "

$ws.Range("D31").Value = "0.0₃0776"
$ws.Range("E31").Value = "  -6.84%  "

$ws.Range("E32").Value = "  -5.40%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.19"
$ws.Range("E33").Value = "  -11.00%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "159.84"
$ws.Range("E34").Value = "  -1.71%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.996"
$ws.Range("E35").Value = "  -0.26%  "

$ws.Range("E36").Value = "  +4.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.53"
$ws.Range("E37").Value = "  -3.06%  "

$ws.Range("E38").Value = "  -10.20%  "

$ws.Range("E39").Value = "  -8.94%  "

$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "303.82"
$ws.Range("E40").Value = "  -6.68%  "

$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.85"
$ws.Range("E41").Value = "  -1.70%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.58"
$ws.Range("E42").Value = "  -7.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.821"
$ws.Range("E43").Value = "  -8.67%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.66"
$ws.Range("E44").Value = "  -6.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.995"
$ws.Range("E45").Value = "  -0.22%  "

$ws.Range("E46").Value = "  -1.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.78"
$ws.Range("E47").Value = "  -1.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.28"
$ws.Range("E48").Value = "  +2.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0928"
$ws.Range("E49").Value = "  -3.79%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.59"
$ws.Range("E50").Value = "  -4.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0516"
$ws.Range("E51").Value = "  -5.26%  "
